$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Daily Orders" ---
$ws = $wb.Worksheets.Item(1)

# Insert a new row at the top of the data (row 2), shifting existing orders down.
$ws.Rows.Item(2).Insert()

# Fill in the details of the new order (was previously row 2, now highest Order ID).
$ws.Cells.Item(2, 1).Value  = 4
$ws.Cells.Item(2, 2).Value  = "2026-01-13 10:51"
$ws.Cells.Item(2, 3).Value  = "Ajay Dwarkunde"
$ws.Cells.Item(2, 4).Value  = "'8087172173"
$ws.Cells.Item(2, 5).Value  = "wakad, pune 411057"
$ws.Cells.Item(2, 6).Value  = "Girl Holding Hands Thali x1, Kalash Haldi Kunku (Golden) x1, Kite Haldi Kunku Set x1"
$ws.Cells.Item(2, 7).Value  = 0
$ws.Cells.Item(2, 8).Value  = "NEW"
$ws.Cells.Item(2, 9).Value  = "PENDING"

# --- Sheet 2: "Summary" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 1).Value = 4
$ws2.Cells.Item(2, 2).Value = 4
